$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Workbook calculation properties: turn off the iterative-calc max-change
#    setting (drops the stray iterateDelta="1E-4" left over on calcPr).
# ---------------------------------------------------------------------------
$wb.Application.Iteration = $false
$wb.Application.MaxChange = 0.001

# ---------------------------------------------------------------------------
# 2. Sheet1 ("Data repo metadata"): the three data rows keep their values,
#    but should no longer carry an explicit row-level style/customFormat
#    flag. Clear each row's format, then restore the per-cell styles by
#    pasting formats from a row that still has them (rotating through so a
#    source row is always available).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).ClearFormats()
$ws1.Range("A1:B1").Copy()
$ws1.Range("A2:B2").PasteSpecial(-4122)

$ws1.Rows.Item(3).ClearFormats()
$ws1.Range("A1:B1").Copy()
$ws1.Range("A3:B3").PasteSpecial(-4122)

$ws1.Rows.Item(1).ClearFormats()
$ws1.Range("A2:B2").Copy()
$ws1.Range("A1:B1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Sheet2 ("Schema repo metadata"): reorder the git-metadata rows to
#    Url / Branch / Revision (was Branch / Revision / Url), clear the
#    row-level customFormat flag the same way as sheet1, add a new blank
#    row 4, move the selection there, and fix up the data-validation
#    ranges so each validation still targets the row holding its field.
# ---------------------------------------------------------------------------

# 3a. Rotate the three rows into their new order, carrying values+formats
#     together (row-level Cut/Copy keeps cell styling attached), using row 5
#     as temporary scratch space for the 3-way rotation.
$ws2.Rows.Item(1).Copy()
$ws2.Rows.Item(5).PasteSpecial(-4104)

$ws2.Rows.Item(3).Copy()
$ws2.Rows.Item(1).PasteSpecial(-4104)

$ws2.Rows.Item(2).Copy()
$ws2.Rows.Item(3).PasteSpecial(-4104)

$ws2.Rows.Item(5).Copy()
$ws2.Rows.Item(2).PasteSpecial(-4104)

$ws2.Rows.Item(5).Delete()

# 3b. Drop the row-level customFormat flag left behind by the row copies,
#     the same rotating clear+paste-format trick used for sheet1.
$ws2.Rows.Item(2).ClearFormats()
$ws2.Range("A1:C1").Copy()
$ws2.Range("A2:C2").PasteSpecial(-4122)

$ws2.Rows.Item(3).ClearFormats()
$ws2.Range("A1:C1").Copy()
$ws2.Range("A3:C3").PasteSpecial(-4122)

$ws2.Rows.Item(1).ClearFormats()
$ws2.Range("A2:C2").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# 3c. Recreate the data validations so each one again targets the row that
#     holds its field (Url -> row1, Branch -> row2, Revision -> row3).
$ws2.Range("B1:C3").Validation.Delete()

$ws2.Range("B1:C1").Validation.Add(6, 2, 8, 255)
$v = $ws2.Range("B1:C1").Validation
$v.ErrorTitle = "Url"
$v.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$v.InputTitle = "Url"
$v.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."

$ws2.Range("B3:C3").Validation.Add(6, 2, 8, 255)
$v = $ws2.Range("B3:C3").Validation
$v.ErrorTitle = "Revision"
$v.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$v.InputTitle = "Revision"
$v.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."

$ws2.Range("B2:C2").Validation.Add(6, 2, 8, 255)
$v = $ws2.Range("B2:C2").Validation
$v.ErrorTitle = "Branch"
$v.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$v.InputTitle = "Branch"
$v.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."

# 3d. Add the new blank row 4 (same height/formatting metadata as the rest)
#     and move the worksheet selection to it.
$ws2.Range("C4").Value = "x"
$ws2.Range("C4").ClearContents()
$ws2.Range("C4").Style = "Normal"
$ws2.Rows.Item(4).RowHeight = 15

$ws2.Activate()
$ws2.Rows.Item(4).Select()
